# "Change File Status read"
# The "Status" column (G) for every course row was stored as the text
# string "false". Convert it to a real boolean FALSE value instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows are 2..63 (row 1 is the header); Status lives in column G (7).
for ($row = 2; $row -le 63; $row++) {
    $ws.Cells.Item($row, 7).Value = $false
}

# Leave the selection where the author ended up after the edit.
$ws.Range("F67").Select()
